# Generate Report for Handoff
# Refresh the "Ready for handoff" / "Handback transform failed" rows so
# their Latest Handoff/Handback timestamp reflects the newest report run.

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 10, 11, 12, 13, 14, 15, 16)

$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 4).Value = "2016-03-25 12:27:18"
}

$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "2016-03-25 12:27:13"
}

$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 5).Value = "2016-03-25 12:27:18"
}
